$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 134, pushing existing rows 134-163 down to 135-164.
$ws.Rows.Item(134).Insert()

# Populate the newly inserted row 134 with the new data record.
$ws.Cells.Item(134, 1).Value = 3
$ws.Cells.Item(134, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(134, 3).Value = "Coquimbo"
$ws.Cells.Item(134, 4).Value = 44637
$ws.Cells.Item(134, 5).Value = 5
$ws.Cells.Item(134, 6).Value = 100112030
$ws.Cells.Item(134, 7).Value = "Poroto granado"
$ws.Cells.Item(134, 8).Value = "Sin especificar"
$ws.Cells.Item(134, 9).Value = "Primera"
$ws.Cells.Item(134, 10).Value = 73
$ws.Cells.Item(134, 11).Value = 22000
$ws.Cells.Item(134, 12).Value = 23000
$ws.Cells.Item(134, 13).Value = 22521
$ws.Cells.Item(134, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(134, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(134, 16).Value = 901
$ws.Cells.Item(134, 17).Value = 25
$ws.Cells.Item(134, 18).Value = "Hortaliza"
